# Auto update: 2025-12-05 02:00:49
#
# Daily refresh of the quantum-computing "DECISION" table:
#   - report date moves from 2025-12-03 to 2025-12-05 (all rows)
#   - the per-ticker metrics (close/RSI/5d return/scores/probabilities)
#     are refreshed with the latest run's numbers
#   - the ticker rows are now emitted in a new order:
#       row2 D-Wave(QBTS) / row3 IonQ(IONQ) / row4 Rigetti(RGTI) / row5 IBM(IBM)
#   - MACRO_SIGNAL flips from "bullish" to "neutral"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163   # Excel.XlPasteType.xlPasteValues

function Set-TextValue($cell, [string]$text) {
    # Assigning a date-shaped literal (e.g. "2025-12-05") directly to
    # .Value triggers Excel's smart text-to-date conversion, turning the
    # cell into a real date serial with a date number format. The source
    # file stores these as plain text, so instead we push the text in via
    # a formula ("=""...""") and immediately collapse it to its cached
    # value with copy / paste-special-values, which keeps the cell a
    # plain string without touching its number format/style.
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial($xlPasteValues) | Out-Null
}

$reportDate = "2025-12-05"
$macroSignal = "⚪ 중립 구간"
$pattern = "Pattern"
$judgment = "⛔ 관망하십시오."

$rows = @(
    @{ Row = 2; Name = "D-Wave Quantum Inc."; Ticker = "QBTS";
       D = 27.95;  E = 64.2;               F = 24.72; G = 50; H = 66; I = 70; J = 80; K = 58.7;
       N = 52.43913937059539 },
    @{ Row = 3; Name = "IonQ, Inc."; Ticker = "IONQ";
       D = 54.46;  E = 64.90000000000001;  F = 16.11; G = 60; H = 60; I = 60; J = 63; K = 57.7;
       N = 52.43913937059539 },
    @{ Row = 4; Name = "Rigetti Computing, Inc."; Ticker = "RGTI";
       D = 28.82;  E = 60.5;               F = 12.71; G = 50; H = 63; I = 63; J = 83; K = 55.9;
       N = 52.43913937059539 },
    @{ Row = 5; Name = "International Business Machines"; Ticker = "IBM";
       D = 308.17; E = 53.1;               F = 1.64;  G = 50; H = 66; I = 56; J = 63; K = 54.7;
       N = 52.43913937059539 }
)

foreach ($r in $rows) {
    $row = $r.Row
    Set-TextValue $ws.Cells.Item($row, 1) $reportDate   # A: 날짜
    $ws.Cells.Item($row, 2).Value = $r.Name              # B: 종목명
    $ws.Cells.Item($row, 3).Value = $r.Ticker            # C: 티커
    $ws.Cells.Item($row, 4).Value = $r.D                 # D: 종가
    $ws.Cells.Item($row, 5).Value = $r.E                 # E: RSI
    $ws.Cells.Item($row, 6).Value = $r.F                 # F: 5일수익률
    $ws.Cells.Item($row, 7).Value = $r.G                 # G: 점수(룰)
    $ws.Cells.Item($row, 8).Value = $r.H                 # H: 3일상승확률(%)
    $ws.Cells.Item($row, 9).Value = $r.I                 # I: 5일상승확률(%)
    $ws.Cells.Item($row, 10).Value = $r.J                # J: 10일상승확률(%)
    $ws.Cells.Item($row, 11).Value = $r.K                # K: 최종점수
    $ws.Cells.Item($row, 12).Value = $pattern            # L: 예측방식
    $ws.Cells.Item($row, 13).Value = $judgment           # M: 판단
    $ws.Cells.Item($row, 14).Value = $r.N                # N: MACRO_SCORE
}

$excel.CutCopyMode = $false

# O: MACRO_SIGNAL - every row shared one string before and still shares
# one after, so write the whole block in a single assignment.
$ws.Range("O2:O5").Value = $macroSignal

Write-Host "Applied 2025-12-05 auto update to DECISION sheet"
